# Reporte de monitores mes de junio
# Updates the "Acciones correctivas" table: rows 11-14 get closed out
# (status -> Cerrado, real-close date filled in) with their corrective
# action text replaced by the real action taken; two previously blank
# rows (15-16) become new open items; the sheet selection left on R14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 11 (id 7) : close out, set real close date, replace the
#     corrective-action placeholder with the actual action taken ---
$ws.Range("F11").Value = "Campaña de publicidad."
$ws.Range("J11").Value = 42558
$ws.Range("K11").Value = "Cerrado"
$ws.Rows(11).Hidden = $true

# --- Row 12 (id 8) ---
$ws.Range("F12").Value = "Campaña de publicidad."
$ws.Range("J12").Value = 42558
$ws.Range("K12").Value = "Cerrado"
$ws.Rows(12).Hidden = $true

# --- Row 13 (id 9) ---
$ws.Range("F13").Value = "Tercera campaña de publicidad para Desarrollo y beneficio del CMMi"
$ws.Range("J13").Value = 42557
$ws.Range("K13").Value = "Cerrado"

# --- Row 14 (id 10) ---
$ws.Range("F14").Value = "Aplicar la medicion resultante en el CMMi para diferenciar tareas de soporte (implementacion y garantia)"
$ws.Range("J14").Value = 42528
$ws.Range("K14").Value = "Cerrado"

# --- Row 15 (id 11) : new open item, reusing the "Desviación ... Mayo"
#     deviation, still awaiting corrective action ---
$ws.Range("D15").Value = "Si"
$ws.Range("E15").Value = "Desviación en ventas planeadas para el mes de Mayo"
$ws.Range("F15").Value = "En espera de acciones correctivas"
$ws.Range("G15").Value = 42557
$ws.Range("H15").Value = "Ricardo Novela"
$ws.Range("I15").Value = 42558
$ws.Range("K15").Value = "Abierto"
$ws.Rows(15).RowHeight = 57

# --- Row 16 (id 12) : new open item, reusing the "Presentar esfuerzos
#     superiores..." deviation, still awaiting corrective action ---
$ws.Range("E16").Value = "Presentar esfuerzos superiores a los estimados en la cotización anual de la empresa lo cual provoca perdidas monetarias a la empresa"
$ws.Range("F16").Value = "En espera de acciones correctivas"
$ws.Range("G16").Value = 42557
$ws.Range("H16").Value = "Ricardo Novela"
$ws.Range("I16").Value = 42558
$ws.Range("K16").Value = "Abierto"
$ws.Rows(16).RowHeight = 114

# --- Leave the sheet's last active selection where the author left it ---
$ws.Range("R14").Select()
